$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.08530170911896
$ws.Range("D2").Value = 1.085124543965623
$ws.Range("E2").Value = 1.08726554938127
$ws.Range("F2").Value = 1.096873305803912
$ws.Range("I2").Value = 1.059245209726623
$ws.Range("J2").Value = 1.090157246685577
$ws.Range("K2").Value = 1.087785785571798
$ws.Range("L2").Value = 1.089921258391643
$ws.Range("M2").Value = 1.099504473980934
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.08701399729553
$ws.Range("D3").Value = 1.086493387049638
$ws.Range("E3").Value = 1.088776999058817
$ws.Range("F3").Value = 1.098384363779025
$ws.Range("I3").Value = 1.059759922991254
$ws.Range("J3").Value = 1.091529202272312
$ws.Range("K3").Value = 1.088972439691654
$ws.Range("L3").Value = 1.091250581336528
$ws.Range("M3").Value = 1.100835199105746
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.088119578485727
$ws.Range("D4").Value = 1.087376886973792
$ws.Range("E4").Value = 1.089752568116499
$ws.Range("F4").Value = 1.099359950300613
$ws.Range("I4").Value = 1.060090394316215
$ws.Range("J4").Value = 1.092414179809803
$ws.Range("K4").Value = 1.089737510452409
$ws.Range("L4").Value = 1.092107789751352
$ws.Range("M4").Value = 1.101693581356674
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.088583806880932
$ws.Range("D5").Value = 1.087747785041193
$ws.Range("E5").Value = 1.090162124136733
$ws.Range("F5").Value = 1.09976957797388
$ws.Range("I5").Value = 1.060228710566107
$ws.Range("J5").Value = 1.092785571623835
$ws.Range("K5").Value = 1.090058490920425
$ws.Range("L5").Value = 1.092467463566305
$ws.Range("M5").Value = 1.102053811967877
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.088661720521894
$ws.Range("D6").Value = 1.087810029923796
$ws.Range("E6").Value = 1.090230857096989
$ws.Range("F6").Value = 1.099838326720327
$ws.Range("I6").Value = 1.060251898600483
$ws.Range("J6").Value = 1.092847891963497
$ws.Range("K6").Value = 1.090112346781093
$ws.Range("L6").Value = 1.092527813828687
$ws.Range("M6").Value = 1.102114259417341
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.088125783703027
$ws.Range("D7").Value = 1.087381844983107
$ws.Range("E7").Value = 1.089758042861338
$ws.Range("F7").Value = 1.09936542575092
$ws.Range("I7").Value = 1.060092244910821
$ws.Range("J7").Value = 1.092419144916044
$ws.Range("K7").Value = 1.089741801971791
$ws.Range("L7").Value = 1.092112598453542
$ws.Range("M7").Value = 1.101698397250678
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.085880886192939
$ws.Range("D8").Value = 1.085587618891338
$ws.Range("E8").Value = 1.087776862679382
$ws.Range("F8").Value = 1.097384430298239
$ws.Range("I8").Value = 1.059419696870241
$ws.Range("J8").Value = 1.090621485142613
$ws.Range("K8").Value = 1.088187400524057
$ws.Range("L8").Value = 1.090371126791929
$ws.Range("M8").Value = 1.099954760770665
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.081906259846069
$ws.Range("D9").Value = 1.082408441323235
$ws.Range("E9").Value = 1.084266591487533
$ws.Range("F9").Value = 1.093876584021225
$ws.Range("I9").Value = 1.058214604362877
$ws.Range("J9").Value = 1.087432105509635
$ws.Range("K9").Value = 1.085426707197249
$ws.Range("L9").Value = 1.087279359456257
$ws.Range("M9").Value = 1.096861242041422
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.079243043975719
$ws.Range("D10").Value = 1.080276619470768
$ws.Range("E10").Value = 1.081912826054109
$ws.Range("F10").Value = 1.091525880624595
$ws.Range("I10").Value = 1.057397499911419
$ws.Range("J10").Value = 1.085290630642365
$ws.Range("K10").Value = 1.083571140111884
$ws.Range("L10").Value = 1.085202023534982
$ws.Range("M10").Value = 1.094784156202508
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.078086468508089
$ws.Range("D11").Value = 1.079350450965541
$ws.Range("E11").Value = 1.080890243180369
$ws.Range("F11").Value = 1.090504973939938
$ws.Range("I11").Value = 1.057040370504788
$ws.Range("J11").Value = 1.08435959254955
$ws.Range("K11").Value = 1.082763951871857
$ws.Range("L11").Value = 1.084298537337754
$ws.Range("M11").Value = 1.093881121233627
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.077656340680966
$ws.Range("D12").Value = 1.079005957199657
$ws.Range("E12").Value = 1.080509888140073
$ws.Range("F12").Value = 1.09012529477788
$ws.Range("I12").Value = 1.056907212825952
$ws.Range("J12").Value = 1.084013185225237
$ws.Range("K12").Value = 1.082463557190281
$ws.Range("L12").Value = 1.083962330947784
$ws.Range("M12").Value = 1.093545134544137
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.077748628533044
$ws.Range("D13").Value = 1.079079873896298
$ws.Range("E13").Value = 1.080591499435731
$ws.Range("F13").Value = 1.090206758674656
$ws.Range("I13").Value = 1.056935798491563
$ws.Range("J13").Value = 1.084087517146521
$ws.Range("K13").Value = 1.082528018793969
$ws.Range("L13").Value = 1.084034476228867
$ws.Range("M13").Value = 1.093617230332158
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.078050924791969
$ws.Range("D14").Value = 1.079321984730383
$ws.Range("E14").Value = 1.08085881364534
$ws.Range("F14").Value = 1.090473599182888
$ws.Range("I14").Value = 1.057029373961025
$ws.Range("J14").Value = 1.084330970283132
$ws.Range("K14").Value = 1.08273913282482
$ws.Range("L14").Value = 1.084270758957761
$ws.Range("M14").Value = 1.093853359948129
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.078237109563503
$ws.Range("D15").Value = 1.079471094111215
$ws.Range("E15").Value = 1.081023445149021
$ws.Range("F15").Value = 1.090637945901021
$ws.Range("I15").Value = 1.057086961958141
$ws.Range("J15").Value = 1.084480892782469
$ws.Range("K15").Value = 1.082869131419192
$ws.Range("L15").Value = 1.084416259153385
$ws.Range("M15").Value = 1.093998772725931
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.079319729427144
$ws.Range("D16").Value = 1.080338020448664
$ws.Range("E16").Value = 1.08198061901387
$ws.Range("F16").Value = 1.091593569776109
$ws.Range("I16").Value = 1.057421131047644
$ws.Range("J16").Value = 1.085352340168406
$ws.Range("K16").Value = 1.083624631335352
$ws.Range("L16").Value = 1.085261899890334
$ws.Range("M16").Value = 1.094844009885401
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.079997910940931
$ws.Range("D17").Value = 1.080880988954002
$ws.Range("E17").Value = 1.082580112711141
$ws.Range("F17").Value = 1.092192185441919
$ws.Range("I17").Value = 1.05762985438964
$ws.Range("J17").Value = 1.085897959130219
$ws.Range("K17").Value = 1.084097534134985
$ws.Range("L17").Value = 1.085791272000037
$ws.Range("M17").Value = 1.095373221059257
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.08039315716033
$ws.Range("D18").Value = 1.081197397330005
$ws.Range("E18").Value = 1.082929461365665
$ws.Range("F18").Value = 1.092541055712005
$ws.Range("I18").Value = 1.057751279437119
$ws.Range("J18").Value = 1.086215846997696
$ws.Range("K18").Value = 1.084373012587571
$ws.Range("L18").Value = 1.086099662053733
$ws.Range("M18").Value = 1.095681550292052
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.080527871151401
$ws.Range("D19").Value = 1.081305234553588
$ws.Range("E19").Value = 1.083048525357449
$ws.Range("F19").Value = 1.092659962293465
$ws.Range("I19").Value = 1.057792628169922
$ws.Range("J19").Value = 1.086324177448962
$ws.Range("K19").Value = 1.084466883266054
$ws.Range("L19").Value = 1.086204750380475
$ws.Range("M19").Value = 1.095786623469054
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.079925182271964
$ws.Range("D20").Value = 1.080822764268151
$ws.Range("E20").Value = 1.082515826514641
$ws.Range("F20").Value = 1.092127989953694
$ws.Range("I20").Value = 1.057607493462699
$ws.Range("J20").Value = 1.085839456959562
$ws.Range("K20").Value = 1.08404683322203
$ws.Range("L20").Value = 1.085734515165348
$ws.Range("M20").Value = 1.095316478068794
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.077961920674696
$ws.Range("D21").Value = 1.079250702289126
$ws.Range("E21").Value = 1.080780110797526
$ws.Range("F21").Value = 1.090395034352762
$ws.Range("I21").Value = 1.057001832272235
$ws.Range("J21").Value = 1.084259295468838
$ws.Range("K21").Value = 1.082676980800104
$ws.Range("L21").Value = 1.084201196561877
$ws.Range("M21").Value = 1.093783841197532
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.076724499485266
$ws.Range("D22").Value = 1.078259538968715
$ws.Range("E22").Value = 1.079685769243856
$ws.Range("F22").Value = 1.089302736697671
$ws.Range("I22").Value = 1.056618110794992
$ws.Range("J22").Value = 1.083262432884126
$ws.Range("K22").Value = 1.081812402541453
$ws.Range("L22").Value = 1.083233594694756
$ws.Range("M22").Value = 1.092816969647702
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.077380773335792
$ws.Range("D23").Value = 1.078785237507554
$ws.Range("E23").Value = 1.080266191947922
$ws.Range("F23").Value = 1.089882046436312
$ws.Range("I23").Value = 1.056821807263202
$ws.Range("J23").Value = 1.083791210611448
$ws.Range("K23").Value = 1.082271048059055
$ws.Range("L23").Value = 1.083746878860738
$ws.Range("M23").Value = 1.093329837871188
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.079958046232838
$ws.Range("D24").Value = 1.080849074411966
$ws.Range("E24").Value = 1.082544875686567
$ws.Range("F24").Value = 1.092156998034187
$ws.Range("I24").Value = 1.05761759839106
$ws.Range("J24").Value = 1.085865892692158
$ws.Range("K24").Value = 1.084069743887603
$ws.Range("L24").Value = 1.085760162321552
$ws.Range("M24").Value = 1.095342118867298
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.08293610987633
$ws.Range("D25").Value = 1.08323247168638
$ws.Range("E25").Value = 1.085176424139053
$ws.Range("F25").Value = 1.094785537499505
$ws.Range("I25").Value = 1.058528545480589
$ws.Range("J25").Value = 1.088259271521975
$ws.Range("K25").Value = 1.086143033851411
$ws.Range("L25").Value = 1.088081455635591
$ws.Range("M25").Value = 1.097663543357359
